$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5372.75
$ws.Range("J64").Value = 5364
$ws.Range("L64").Value = 5364
$ws.Range("N64").Value = -5860
$ws.Range("H67").Value = 5372.75
$ws.Range("J67").Value = 5364
$ws.Range("L67").Value = 5364
$ws.Range("N67").Value = -7080
$ws.Range("H74").Value = 30581.223
$ws.Range("I74").Value = 51804.75
$ws.Range("J74").Value = 13602.4
$ws.Range("K74").Value = 51804.75
$ws.Range("L74").Value = 13602.4
$ws.Range("M74").Value = -50868.75
$ws.Range("N74").Value = -15474.4
$ws.Range("H77").Value = 30581.223
$ws.Range("I77").Value = 51804.75
$ws.Range("J77").Value = 13602.4
$ws.Range("K77").Value = 259023.75
$ws.Range("L77").Value = 68012
$ws.Range("M77").Value = -254343.75
$ws.Range("N77").Value = -77372
$ws.Range("H88").Value = 7043.6665
$ws.Range("I88").Value = 1649
$ws.Range("J88").Value = 8585
$ws.Range("K88").Value = 1649
$ws.Range("L88").Value = 8585
$ws.Range("M88").Value = -1243
$ws.Range("N88").Value = -9397
$ws.Range("H91").Value = 7043.6665
$ws.Range("I91").Value = 1649
$ws.Range("J91").Value = 8585
$ws.Range("K91").Value = 1649
$ws.Range("L91").Value = 8585
$ws.Range("M91").Value = -245
$ws.Range("N91").Value = -11393
$ws.Range("H92").Value = 169.27272
$ws.Range("I92").Value = 125.28571
$ws.Range("K92").Value = 125.28571
$ws.Range("M92").Value = 1122.71429
$ws.Range("H98").Value = 965.6667
$ws.Range("I98").Value = 965.6667
$ws.Range("K98").Value = 965.6667
$ws.Range("M98").Value = 532.3333
$ws.Range("H122").Value = 965.6667
$ws.Range("I122").Value = 965.6667
$ws.Range("K122").Value = 2897.0001
$ws.Range("M122").Value = -447.0001000000002
$ws.Range("H135").Value = 570.1667
$ws.Range("I135").Value = 570.1667
$ws.Range("K135").Value = 5131.5003
$ws.Range("M135").Value = -2596.5003
$ws.Range("H138").Value = 4124.2285
$ws.Range("J138").Value = 4966.4585
$ws.Range("L138").Value = 14899.3755
$ws.Range("N138").Value = -25179.3755

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 366.4
$ws.Range("I2").Value = 366.4
$ws.Range("K2").Value = 366.4
$ws.Range("M2").Value = -253.4
$ws.Range("H5").Value = 148.28572
$ws.Range("I5").Value = 148.28572
$ws.Range("K5").Value = 148.28572
$ws.Range("M5").Value = -36.28572
$ws.Range("H32").Value = 2408855.2
$ws.Range("I32").Value = 3184693.5
$ws.Range("J32").Value = 702010.9
$ws.Range("K32").Value = 3184693.5
$ws.Range("L32").Value = 702010.9
$ws.Range("M32").Value = -3184406.5
$ws.Range("N32").Value = -702584.9
$ws.Range("H88").Value = 2023.8
$ws.Range("I88").Value = 1718
$ws.Range("J88").Value = 2635.4
$ws.Range("K88").Value = 1718
$ws.Range("L88").Value = 2635.4
$ws.Range("M88").Value = -1312
$ws.Range("N88").Value = -3447.4
$ws.Range("H91").Value = 2023.8
$ws.Range("I91").Value = 1718
$ws.Range("J91").Value = 2635.4
$ws.Range("K91").Value = 1718
$ws.Range("L91").Value = 2635.4
$ws.Range("M91").Value = -314
$ws.Range("N91").Value = -5443.4
$ws.Range("H116").Value = 366.4
$ws.Range("I116").Value = 366.4
$ws.Range("K116").Value = 366.4
$ws.Range("M116").Value = 1927.6
$ws.Range("H122").Value = 2352
$ws.Range("I122").Value = 1808
$ws.Range("K122").Value = 5424
$ws.Range("M122").Value = -2974
$ws.Range("H132").Value = 2697.0908
$ws.Range("I132").Value = 2697.0908
$ws.Range("K132").Value = 8091.2724
$ws.Range("M132").Value = -5561.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 366.4
$ws.Range("I3").Value = 366.4
$ws.Range("K3").Value = 366.4
$ws.Range("M3").Value = -252.4
$ws.Range("H4").Value = 148.28572
$ws.Range("I4").Value = 148.28572
$ws.Range("K4").Value = 148.28572
$ws.Range("M4").Value = -33.28572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2708.25
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null
$ws.Range("H126").Value = 2708.25
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 939.4
$ws.Range("I5").Value = 299.25
$ws.Range("K5").Value = 897.75
$ws.Range("M5").Value = -785.75
$ws.Range("H121").Value = 3808.8572
$ws.Range("I121").Value = 1406.75
$ws.Range("K121").Value = 4220.25
$ws.Range("M121").Value = -2910.25
$ws.Range("H135").Value = 939.4
$ws.Range("I135").Value = 299.25
$ws.Range("K135").Value = 2693.25
$ws.Range("M135").Value = -158.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2854.8572
$ws.Range("I102").Value = 2854.8572
$ws.Range("K102").Value = 2854.8572
$ws.Range("M102").Value = -1232.8572
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 474.85715
$ws.Range("I16").Value = 474.85715
$ws.Range("K16").Value = 474.85715
$ws.Range("M16").Value = -304.85715
$ws.Range("H46").Value = 3916.9
$ws.Range("I46").Value = 3710.2856
$ws.Range("K46").Value = 3710.2856
$ws.Range("M46").Value = -3522.2856
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256
$ws.Range("H122").Value = 6559
$ws.Range("J122").Value = 7641.143
$ws.Range("L122").Value = 22923.429
$ws.Range("N122").Value = -27823.429
$ws.Range("H132").Value = 6883.231
$ws.Range("I132").Value = 7040.0835
$ws.Range("K132").Value = 21120.2505
$ws.Range("M132").Value = -18590.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 4199
$ws.Range("J19").Value = 4199
$ws.Range("L19").Value = 4199
$ws.Range("N19").Value = -4547
$ws.Range("H75").Value = 75084.25
$ws.Range("I75").Value = 73108
$ws.Range("J75").Value = 75743
$ws.Range("K75").Value = 73108
$ws.Range("L75").Value = 75743
$ws.Range("M75").Value = -72172
$ws.Range("N75").Value = -77615
$ws.Range("H78").Value = 75084.25
$ws.Range("I78").Value = 73108
$ws.Range("J78").Value = 75743
$ws.Range("K78").Value = 219324
$ws.Range("L78").Value = 227229
$ws.Range("M78").Value = -214644
$ws.Range("N78").Value = -236589
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 5475
$ws.Range("I132").Value = 5475
$ws.Range("K132").Value = 16425
$ws.Range("M132").Value = -13895
